# Apply the "Deploying to gh-pages ... LinuxForHealth" update to the
# StructureDefinition-documented-system workbook.
#
# Sheet "Metadata": rebrand from ibm.com/Alvearie to linuxforhealth.org,
#   bump Version 7.0.0 -> 8.0.0, and refresh the Date.
# Sheet "Elements": the Fixed Value for Extension.url (same URL as above)
#   and the Binding Value Set URL also move from ibm.com to
#   linuxforhealth.org; also the ele-1/ext-1 Constraint(s) note that was
#   duplicated onto the top-level "Extension" row is cleared there (it
#   correctly stays on the "Extension.extension" row only).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/documented-system"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/documented-system"
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/legal-document-system"
$elements.Range("AI2").Value = ""
